$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 39 (shifts old rows 39-43 down to 40-44) ---
$ws.Rows("39:39").Insert()

# --- Row 38: E38 becomes a number (2), F38 becomes a date (43903) with date formatting ---
$ws.Range("F35").Copy()
$ws.Range("F38").PasteSpecial(-4122)  # xlPasteFormats (copy date number format/style from F35)
$ws.Range("E38").Value = 2
$ws.Range("F38").Value = 43903

# --- Row 39 (new row): "Fix Incompatible Code In Admin Functionality - Technical" task ---
$ws.Range("A39").Value = "Fix Incompatible Code In Admin Functionality - Technical"
$ws.Range("B39").Value = "Admin must have all necessary functionalities required"
$ws.Range("C39").Value = "Gentian Gashi"
$ws.Range("D39").Value = 3
$ws.Range("E39").Value = 3
$ws.Range("F39").Value = 43903

Write-Output "done"
